$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.517.59"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "1.819.92"
$ws.Range("E3").Value = "  -0.60%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.33"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5157"
$ws.Range("E7").Value = "  -3.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3874"
$ws.Range("E8").Value = "  -2.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08448"
$ws.Range("E9").Value = "  +8.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.92"
$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.405"
$ws.Range("E12").Value = "  +0.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.01"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.002"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.499"
$ws.Range("E15").Value = "  -1.07%  "

$ws.Range("D16").Value = "1.812.81"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001138"
$ws.Range("E17").Value = "  +4.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.69"
$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06695"
$ws.Range("E19").Value = "  +1.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.74"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.081"
$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("D23").Value = "28.539.08"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +2.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.273"
$ws.Range("E25").Value = "  +1.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.12"
$ws.Range("E26").Value = "  +1.68%  "

$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.424"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "2.023.81"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.05"
$ws.Range("E30").Value = "  +0.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.096"
$ws.Range("E31").Value = "  -4.15%  "

$ws.Range("E32").Value = "  -3.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.755"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07482"
$ws.Range("E34").Value = "  +2.59%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.684"
$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2226"
$ws.Range("E36").Value = "  -1.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02364"
$ws.Range("E37").Value = "  +0.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.213"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.751"
$ws.Range("E39").Value = "  -2.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6333"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.24"
$ws.Range("E41").Value = "  -1.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.190"
$ws.Range("E42").Value = "  -0.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.400"
$ws.Range("E43").Value = "  +0.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.60"
$ws.Range("E44").Value = "  +0.26%  "

$ws.Range("E45").Value = "  +1.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5908"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.02"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.989"
$ws.Range("E48").Value = "  -0.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.198"
$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06982"
$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.17"
$ws.Range("E51").Value = "  -0.88%  "
